$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header labels: "_old" -> "_FV2404", "_new" -> "_FV2410"
$oldHeaders = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
$newHeaders = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}
for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# Create an Excel Table (ListObject) over A1:U92 with autofilter
$range = $ws.Range("A1:U92")
$table = $ws.ListObjects.Add(1, $range, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (split/pane)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
